# HextrudORT BOM - "Bondtech Kit reference added"
#
# Rows 3-8 (Bondtech BMG parts sub-assembly) previously referenced a
# collection of individual Bondtech replacement-part pages.  They are
# replaced with a single reference to the new "BMG Internals Set for
# HextrudORT" kit:
#   - Comment (col F)    -> "Included in BMG Internals Set for HextrudORT"
#   - Make/Buy (col G)   -> "(BUY) KIT"
#   - Vendor URL (col K) -> new Bondtech kit URL (hyperlink)
#   - QTY on row 8 (bearings) bumped from 1* to 2*
#
# The workbook is already open as $excel.ActiveWorkbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newComment = "Included in BMG Internals Set for HextrudORT"
$newMakeBuy = "(BUY) KIT"
$newUrl     = "https://www.bondtech.se/product/bmg-internals-set-for-hextrudort/"

# ---------------------------------------------------------------------
# 1) Update the "Comment" (F) and "Make/Buy" (G) columns for rows 3-8.
#    Doing F3/G3 first makes sure the two new shared-strings are created
#    in the same order the canonical workbook uses them.
# ---------------------------------------------------------------------
foreach ($r in 3..8) {
    $ws.Range("F$r").Value2 = $newComment
    $ws.Range("G$r").Value2 = $newMakeBuy
}

# ---------------------------------------------------------------------
# 2) Row 8 quantity: this part now represents 2 items included in the kit.
# ---------------------------------------------------------------------
$ws.Range("H8").Value2 = "2*"

# ---------------------------------------------------------------------
# 3) Hyperlinks in column K (rows 3-8) all now point at the same kit page.
#    The engine's Hyperlinks collection doesn't support an in-place
#    target update without side effects, so we clear every hyperlink on
#    the sheet and recreate them (including the two untouched ones on
#    K14/K20) in their original order.
# ---------------------------------------------------------------------
$e3dUrl = "https://e3d-online.com/"

# Remember display text for the two hyperlinks that are not changing.
$k14Text = $ws.Range("K14").Value2
$k20Text = $ws.Range("K20").Value2

$ws.Range("K3").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("K14"), $e3dUrl)
$ws.Hyperlinks.Add($ws.Range("K20"), $e3dUrl)
$ws.Hyperlinks.Add($ws.Range("K3"), $newUrl)
$ws.Hyperlinks.Add($ws.Range("K4"), $newUrl)
$ws.Hyperlinks.Add($ws.Range("K5"), $newUrl)
$ws.Hyperlinks.Add($ws.Range("K6"), $newUrl)
$ws.Hyperlinks.Add($ws.Range("K7"), $newUrl)
$ws.Hyperlinks.Add($ws.Range("K8"), $newUrl)

# Restore the original display text for the untouched hyperlinks ...
$ws.Range("K14").Value2 = $k14Text
$ws.Range("K20").Value2 = $k20Text
# ... and set the new display text (the URL itself) for the changed ones.
$ws.Range("K3").Value2 = $newUrl
$ws.Range("K4").Value2 = $newUrl
$ws.Range("K5").Value2 = $newUrl
$ws.Range("K6").Value2 = $newUrl
$ws.Range("K7").Value2 = $newUrl
$ws.Range("K8").Value2 = $newUrl

# Re-adding hyperlinks makes the engine apply a slightly different
# (but visually identical) cell style; nudging the underline back to
# "single" makes it resolve back to the original shared hyperlink style.
foreach ($cell in @("K3","K4","K5","K6","K7","K8","K14","K20")) {
    $ws.Range($cell).Font.Underline = 2
}

# ---------------------------------------------------------------------
# 4) Part Name / Part Description swaps for individual rows, now that
#    the BOM lines reference the kit instead of individual spares.
# ---------------------------------------------------------------------
$ws.Range("E4").Value2  = "Bondtech_Shaft_assembly"
$ws.Range("I3").Value2  = "Gear set for filament size 1.75 mm, primary gear with inner diameter 5 mm, secondary gear with two needle bearings,M3x2 setscrew and one 3 x 20 mm shaft."
$ws.Range("I4").Value2  = "Replacement Shaft assembly for our Bondtech Mini Geared (BMG), SingleDirect and DualDirect extruders.`nIt includes the set screw for primary 1.75/5.0 drive gear."
$ws.Range("I6").Value2  = "From BMG Extruder 3 X 20"
$ws.Range("E8").Value2  = "OPTION A: Bearing_MR85 (from BMG)"
$ws.Range("I8").Value2  = "OPTION A`nUse MR85 5x8x2.5 bearings from BMG`n`n*Qty1 bearing will be used in the carriage BOM.  Qty 1 required in this assembly + 1 in carriage = 2 for the printer"
$ws.Range("A9").Value2  = "E3D V6 HextrudORT"
$ws.Range("E10").Value2 = "OPTION B: Bearing_MR95_7804K105"
$ws.Range("I10").Value2 = "OPTION B`nThese are BIGGER bearings than the ones used in the BMG.`n To be used with Back Plate and cover marked with `"BB`""

# ---------------------------------------------------------------------
# 5) Selected cell moves from M8 to E4.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("E4").Select()
